$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix formatting on A38 so it matches the rest of the data rows (style was
# previously an outlier with no explicit font/format application).
$ws.Range("A38").Style = "Normal"

# New session data: two additional climbing-log rows (39 and 40).
$ws.Range("A39").Value = "3 Mar 2021"
$ws.Range("B39").Value = "0,6"
$ws.Range("C39").Value = "-12.5,6"
$ws.Range("D39").Value = "-22.5,6"
$ws.Range("E39").Value = "-42.5,6"
$ws.Range("F39").Value = "-30,5,9"
$ws.Range("G39").Value = "-15,5,9"
$ws.Range("H39").Value = "-40,4,6,5"
$ws.Range("I39").Value = "-22.5,4,8,7"
$ws.Range("J39").Value = "-30,4,8,7"

$ws.Range("A40").Value = "10 Mar 2021"
$ws.Range("B40").Value = "0,6"
$ws.Range("C40").Value = "-10,6"
$ws.Range("D40").Value = "-20,6"
$ws.Range("E40").Value = "-40,5,5"
$ws.Range("F40").Value = "-22.5,3,7,7,3"
$ws.Range("G40").Value = "-12.5,6"
$ws.Range("H40").Value = "-37.5,4,9,5"
$ws.Range("I40").Value = "-20,5,7"
$ws.Range("J40").Value = "-27.5,6"

# Make the new rows share the same cell formatting as the rest of the table.
$ws.Range("A39:J40").Style = "Normal"

# Move the active selection to the new last cell, as it is after manual entry.
[void]$ws.Range("J40").Select()
